$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "metadata" sheet: fix the date-range text, add a per-row "note" column
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("metadata")

# Fix stray spaces around the dash in the tx.change period value.
$meta.Range("C2").Value = "1960-2022"

# Insert a brand-new column D ("note"); this pushes the old, wide D column
# (the long explanatory note, with its formatting) one slot right to E.
$meta.Columns.Item(4).Insert()

# Populate the new "note" column.
$meta.Range("D1").Value = "note"
$meta.Range("D2").Value = "une moyenne annuelle sur la base des moyennes mensuelles."
$meta.Range("D3").Value = "n/a"
$meta.Range("D4").Value = "n/a"
$meta.Range("D5").Value = "n/a"
$meta.Range("D6").Value = "n/a"
$meta.Range("D7").Value = "n/a"

# The long note text that used to live in D now moves to its own sheet, so
# clear the text that slid into column E but keep that cell's formatting.
$meta.Range("E1").ClearContents()
$meta.Range("E2").ClearContents()

# Match the target column widths.
$meta.Columns.Item(4).ColumnWidth = 101.5
$meta.Columns.Item(5).ColumnWidth = 225.6640625

# ---------------------------------------------------------------------------
# 2) Add a new sheet with the full explanatory note for tx.change
# ---------------------------------------------------------------------------
$explic = $wb.Worksheets.Add($null, $meta)
$explic.Name = "explication.tx.change"

$explic.Range("A1").Value = "tx.change"
$explic.Range("A2").Value = "Le taux de change officiel fait référence au taux de change indiqué par les autorités nationales ou au taux fixé légalement par le marché des taux. Il est calculé comme étant une moyenne annuelle sur la base des moyennes mensuelles (unités de devises locales par rapport au dollar américain)."
$explic.Range("A2").Font.Size = 11
$explic.Range("A2").Font.Color = 0

$explic.Columns.Item(1).ColumnWidth = 225.6640625

$explic.Activate()
$explic.Range("A2").Select()

# ---------------------------------------------------------------------------
# 3) Leave "metadata" as the active / selected sheet and cell, matching the
#    target workbook view state.
# ---------------------------------------------------------------------------
$meta.Activate()
$meta.Range("B7").Select()
